# Daily attendance processing - 2025-10-12 21:40:39
#
# The "Recorded By" column (G) lists the users who recorded each
# attendance session, separated by ", ". A couple of entries had the
# "System" user listed first/out of order relative to the human
# recorder; this pass normalizes those entries by moving the
# "System"/"system" token(s) to the end of the list.
#
# Concretely:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system"   -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "Recorded By" values. Use Find/FindNext so that only
# cells that actually contain one of the two known (unordered) strings are
# touched - cells elsewhere in the column are left completely untouched.
$col = $ws.Columns.Item(7)

function Swap-RecordedBy($oldText, $newText) {
    $first = $col.Find($oldText)
    if (-not $first) {
        return
    }

    $firstAddr = $first.Address()
    $current = $first
    $visited = @{}

    while ($true) {
        $addr = $current.Address()
        if ($visited.ContainsKey($addr)) {
            break
        }
        $visited[$addr] = $true

        if ($current.Text -eq $oldText) {
            $current.Value = $newText
        }

        $next = $col.FindNext($current)
        if ((-not $next) -or ($next.Address() -eq $firstAddr)) {
            break
        }
        $current = $next
    }
}

Swap-RecordedBy "System, dnasr281@gmail.com" "dnasr281@gmail.com, System"
Swap-RecordedBy "backup@backdoor.com, System, system" "backup@backdoor.com, system, System"
